$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on columns D and E so numeric-looking strings are not
# auto-converted to numbers by Excel (values must stay text, same as source).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "35.249.36"
$ws.Range("E2").Value = "  +1.38%  "

$ws.Range("D3").Value = "1.889.28"
$ws.Range("E3").Value = "  +1.33%  "

$ws.Range("D5").Value = "246.05"
$ws.Range("E5").Value = "  +0.55%  "

$ws.Range("D6").Value = "0.686"
$ws.Range("E6").Value = "  +1.44%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "42.81"
$ws.Range("E8").Value = "  +3.54%  "

$ws.Range("D9").Value = "0.356"
$ws.Range("E9").Value = "  +4.19%  "

$ws.Range("D10").Value = "54.51"
$ws.Range("E10").Value = "  +6.45%  "

$ws.Range("D11").Value = "0.0742"
$ws.Range("E11").Value = "  +2.06%  "

$ws.Range("D12").Value = "0.0977"
$ws.Range("E12").Value = "  +1.26%  "

$ws.Range("D13").Value = "13.83"
$ws.Range("E13").Value = "  +8.30%  "

$ws.Range("D14").Value = "0.784"
$ws.Range("E14").Value = "  +11.12%  "

$ws.Range("D15").Value = "2.161.35"
$ws.Range("E15").Value = "  +1.28%  "

$ws.Range("E16").Value = "  +3.15%  "

$ws.Range("D17").Value = "1.885.34"
$ws.Range("E17").Value = "  +0.39%  "

$ws.Range("D18").Value = "35.272.38"
$ws.Range("E18").Value = "  +1.55%  "

$ws.Range("D19").Value = "73.29"
$ws.Range("E19").Value = "  +1.82%  "

$ws.Range("D20").Value = "0.0₃0823"
$ws.Range("E20").Value = "  +2.00%  "

$ws.Range("D21").Value = "243.41"
$ws.Range("E21").Value = "  +0.72%  "

$ws.Range("D22").Value = "12.77"
$ws.Range("E22").Value = "  +2.06%  "

$ws.Range("D23").Value = "5.18"
$ws.Range("E23").Value = "  +6.57%  "

$ws.Range("D24").Value = "2.64"
$ws.Range("E24").Value = "  +6.88%  "

$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("D26").Value = "2.15"
$ws.Range("E26").Value = "  -0.89%  "

$ws.Range("D27").Value = "167.06"
$ws.Range("E27").Value = "  +1.98%  "

$ws.Range("D28").Value = "8.46"
$ws.Range("E28").Value = "  +1.58%  "

$ws.Range("D29").Value = "18.25"
$ws.Range("E29").Value = "  +1.24%  "

$ws.Range("D30").Value = "0.127"
$ws.Range("E30").Value = "  +1.02%  "

$ws.Range("D31").Value = "4.34"
$ws.Range("E31").Value = "  +4.12%  "

$ws.Range("D32").Value = "0.0594"
$ws.Range("E32").Value = "  +3.61%  "

$ws.Range("D33").Value = "4.18"
$ws.Range("E33").Value = "  +2.52%  "

$ws.Range("E34").Value = "  +15.77%  "

$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("E36").Value = "  -13.63%  "

$ws.Range("D37").Value = "0.849"
$ws.Range("E37").Value = "  +3.73%  "

$ws.Range("D38").Value = "1.94"
$ws.Range("E38").Value = "  -0.28%  "

$ws.Range("D39").Value = "0.0712"
$ws.Range("E39").Value = "  +6.73%  "

$ws.Range("E40").Value = "  +5.25%  "

$ws.Range("D41").Value = "98.16"
$ws.Range("E41").Value = "  +0.76%  "

$ws.Range("D42").Value = "17.01"
$ws.Range("E42").Value = "  +0.43%  "

$ws.Range("E43").Value = "  +0.64%  "

$ws.Range("D44").Value = "1.328.23"
$ws.Range("E44").Value = "  +4.07%  "

$ws.Range("D45").Value = "13.49"
$ws.Range("E45").Value = "  +14.59%  "

$ws.Range("D46").Value = "2.36"
$ws.Range("E46").Value = "  +3.05%  "

$ws.Range("E47").Value = "  -0.78%  "

$ws.Range("E48").Value = "  +0.60%  "

$ws.Range("D49").Value = "2.73"
$ws.Range("E49").Value = "  +0.60%  "

$ws.Range("D50").Value = "6.26"
$ws.Range("E50").Value = "  +0.18%  "

$ws.Range("D51").Value = "2.062.66"
$ws.Range("E51").Value = "  +1.00%  "
